$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the sample rows 4-7, leaving only the (now-empty) column G cells
# behind as placeholders.
$ws.Range("A4:F7").ClearContents()
$ws.Range("G4:G7").ClearContents()

# Remove bold from the "Pattern?" (G) column's data cells (rows 1-7).
# The column itself (col G) keeps its bold column-style, but each cell
# now carries its own (non-bold) direct formatting, same as when a user
# selects G1:G7 and presses Ctrl+B to toggle bold off.
$ws.Range("G1:G7").Font.Bold = $false

# Reproduce the final selection left by the user: whole rows 4:7 selected
# with the active cell on A4.
$ws.Range("A4:XFD7").Select()
